# FA140_TestData_CapitalizeReverseCapitalizedCIPAssets_21C.xlsx
#
# The author re-uploaded this Selenium/UI-automation test-data workbook and,
# in doing so, scrubbed the hard-coded Oracle Cloud login that had been sitting
# in row 2 of the Input_Value sheet (columns N:P -> URL / UserName / Password).
# The SelectBook/AssetType/AssetNumber/Date values (I2:L2) are left intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# Wipe the stored Oracle Cloud credentials (URL, UserName, Password).
$ws.Range("N2:P2").ClearContents()

# Leave the selection on the cells that were just cleared.
$ws.Range("N2:P2").Select()
